$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and D contain text that Excel would otherwise auto-convert
# (A looks like a date, D looks like a plain integer). Force text entry,
# then clear the format so the cell keeps the default (unstyled) style.
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "2023-06-19"
$ws.Range("A60").ClearFormats()

$ws.Range("B60").Value = "22:10:23"
$ws.Range("C60").Value = "Monday"

$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "25"
$ws.Range("D60").ClearFormats()

$ws.Range("E60").Value = 122028
$ws.Range("F60").Value = 133678
$ws.Range("G60").Value = 162436
$ws.Range("H60").Value = 133092
$ws.Range("I60").Value = 177234
$ws.Range("J60").Value = 114659
$ws.Range("K60").Value = 201549
$ws.Range("L60").Value = 225244
$ws.Range("M60").Value = 175511
$ws.Range("N60").Value = 103775
$ws.Range("O60").Value = 39204
$ws.Range("P60").Value = 33957
$ws.Range("Q60").Value = 51812
$ws.Range("R60").Value = -1
$ws.Range("S60").Value = 36484
$ws.Range("T60").Value = -1
